# fix NPC HP error
# Insert a new "MAXHP" column before the existing "MAXMP" column (G) on Sheet1,
# shifting all subsequent columns one to the right, and seed the new MAXHP
# values from the (post-shift) MAXMP column so the two stay in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember column F's width so the newly inserted column G can match it
# (mirrors the look of the original MAXMP column it is displacing).
$fWidth = $ws.Columns("F").ColumnWidth

# Insert a new blank column at G; everything from G onward (MAXMP..AtkDis)
# shifts right by one (H..Z).
$ws.Columns("G").Insert()
$ws.Columns("G").ColumnWidth = $fWidth

# Header for the newly inserted column.
$ws.Range("G1").Value = "MAXHP"

# Populate MAXHP for each data row with the same value as MAXMP (now in
# column H after the shift) for rows 2-6.
for ($r = 2; $r -le 6; $r++) {
    $maxmp = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 7).Value2 = $maxmp
}

# Restore a sane view: no frozen/scrolled top-left override, selection on I8.
[void]$ws.Range("I8").Select()
